# Add a new slide (3rd slide) using the "Title and Content" layout
# (ppLayoutText = 2 -> slideLayout2.xml "Título y objetos"), matching
# the new ppt/slides/slide3.xml / sldId 261 added in the target deck.
$p = $ppt.ActivePresentation
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Content placeholder (shape 2 = "Marcador de contenido 2"): three runs
# of Spanish ("es-ES_tradnl") text reading "PERSONA(nif, nombre, ".
$content = $s.Shapes.Item(2)
$tr = $content.TextFrame.TextRange
$tr.Text = "PERSONA("
$tr.LanguageID = "es-ES_tradnl"
$tr.InsertAfter("nif")
$tr.InsertAfter(", nombre, ")
